$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('I26').Value = 'sv'
$ws.Range('J26').Value = 'Statement-opinion'
$ws.Range('I29').Value = 'sv'
$ws.Range('J29').Value = 'Statement-opinion'
$ws.Range('I30').Value = 'sv'
$ws.Range('J30').Value = 'Statement-opinion'
$ws.Range('I31').Value = 'b'
$ws.Range('J31').Value = 'Acknowledge (Backchannel)'
$ws.Range('I34').Value = 'sd'
$ws.Range('J34').Value = 'Statement-non-opinion'
$ws.Range('I45').Value = 'sd'
$ws.Range('J45').Value = 'Statement-non-opinion'
$ws.Range('I46').Value = 'sv'
$ws.Range('J46').Value = 'Statement-opinion'
$ws.Range('I93').Value = 'sv'
$ws.Range('J93').Value = 'Statement-opinion'
$ws.Range('I94').Value = 'sv'
$ws.Range('J94').Value = 'Statement-opinion'
$ws.Range('I98').Value = 'sv'
$ws.Range('J98').Value = 'Statement-opinion'
$ws.Range('I111').Value = 'sv'
$ws.Range('J111').Value = 'Statement-opinion'
$ws.Range('I136').Value = 'sd'
$ws.Range('J136').Value = 'Statement-non-opinion'
$ws.Range('I153').Value = 'ba'
$ws.Range('J153').Value = 'Appreciation'
$ws.Range('I170').Value = 'ba'
$ws.Range('J170').Value = 'Appreciation'
$ws.Range('I175').Value = '%'
$ws.Range('J175').Value = 'Uninterpretable'
$ws.Range('I184').Value = 'sd'
$ws.Range('J184').Value = 'Statement-non-opinion'
$ws.Range('I185').Value = 'sd'
$ws.Range('J185').Value = 'Statement-non-opinion'
$ws.Range('I186').Value = 'b'
$ws.Range('J186').Value = 'Acknowledge (Backchannel)'
$ws.Range('I207').Value = 'aa'
$ws.Range('J207').Value = 'Agree/Accept'
$ws.Range('I234').Value = 'aa'
$ws.Range('J234').Value = 'Agree/Accept'
$ws.Range('I239').Value = 'sv'
$ws.Range('J239').Value = 'Statement-opinion'
$ws.Range('I252').Value = '%'
$ws.Range('J252').Value = 'Uninterpretable'
$ws.Range('I264').Value = '%'
$ws.Range('J264').Value = 'Uninterpretable'
$ws.Range('I275').Value = 'sv'
$ws.Range('J275').Value = 'Statement-opinion'
$ws.Range('I289').Value = 'sv'
$ws.Range('J289').Value = 'Statement-opinion'
$ws.Range('I290').Value = '%'
$ws.Range('J290').Value = 'Uninterpretable'
$ws.Range('I296').Value = 'b'
$ws.Range('J296').Value = 'Acknowledge (Backchannel)'
$ws.Range('I299').Value = 'sv'
$ws.Range('J299').Value = 'Statement-opinion'
$ws.Range('I312').Value = 'sd'
$ws.Range('J312').Value = 'Statement-non-opinion'
$ws.Range('I317').Value = 'sd'
$ws.Range('J317').Value = 'Statement-non-opinion'
$ws.Range('I325').Value = 'b'
$ws.Range('J325').Value = 'Acknowledge (Backchannel)'
$ws.Range('I326').Value = 'sd'
$ws.Range('J326').Value = 'Statement-non-opinion'
$ws.Range('I333').Value = 'sd'
$ws.Range('J333').Value = 'Statement-non-opinion'
$ws.Range('I334').Value = 'ba'
$ws.Range('J334').Value = 'Appreciation'
$ws.Range('I337').Value = 'aa'
$ws.Range('J337').Value = 'Agree/Accept'
$ws.Range('I338').Value = 'ba'
$ws.Range('J338').Value = 'Appreciation'
$ws.Range('I339').Value = 'b'
$ws.Range('J339').Value = 'Acknowledge (Backchannel)'
$ws.Range('I363').Value = 'aa'
$ws.Range('J363').Value = 'Agree/Accept'
$ws.Range('I364').Value = 'sv'
$ws.Range('J364').Value = 'Statement-opinion'
$ws.Range('I365').Value = 'sd'
$ws.Range('J365').Value = 'Statement-non-opinion'
$ws.Range('I374').Value = '%'
$ws.Range('J374').Value = 'Uninterpretable'
$ws.Range('I379').Value = 'ba'
$ws.Range('J379').Value = 'Appreciation'
$ws.Range('I387').Value = 'aa'
$ws.Range('J387').Value = 'Agree/Accept'
$ws.Range('I408').Value = 'sd'
$ws.Range('J408').Value = 'Statement-non-opinion'
$ws.Range('I412').Value = 'b'
$ws.Range('J412').Value = 'Acknowledge (Backchannel)'
$ws.Range('I432').Value = 'ba'
$ws.Range('J432').Value = 'Appreciation'
$ws.Range('I434').Value = 'sd'
$ws.Range('J434').Value = 'Statement-non-opinion'
$ws.Range('I436').Value = 'aa'
$ws.Range('J436').Value = 'Agree/Accept'
$ws.Range('I441').Value = 'aa'
$ws.Range('J441').Value = 'Agree/Accept'
$ws.Range('I458').Value = 'sd'
$ws.Range('J458').Value = 'Statement-non-opinion'
$ws.Range('I460').Value = 'sv'
$ws.Range('J460').Value = 'Statement-opinion'
$ws.Range('I466').Value = 'sv'
$ws.Range('J466').Value = 'Statement-opinion'
$ws.Range('I467').Value = 'sv'
$ws.Range('J467').Value = 'Statement-opinion'
$ws.Range('I468').Value = 'ba'
$ws.Range('J468').Value = 'Appreciation'
$ws.Range('I470').Value = 'aa'
$ws.Range('J470').Value = 'Agree/Accept'
$ws.Range('I476').Value = 'sv'
$ws.Range('J476').Value = 'Statement-opinion'
$ws.Range('I479').Value = 'sd'
$ws.Range('J479').Value = 'Statement-non-opinion'
$ws.Range('I482').Value = 'sd'
$ws.Range('J482').Value = 'Statement-non-opinion'
$ws.Range('I484').Value = 'sd'
$ws.Range('J484').Value = 'Statement-non-opinion'
$ws.Range('I486').Value = 'ba'
$ws.Range('J486').Value = 'Appreciation'
$ws.Range('I491').Value = 'ba'
$ws.Range('J491').Value = 'Appreciation'
$ws.Range('I502').Value = 'ba'
$ws.Range('J502').Value = 'Appreciation'
$ws.Range('I503').Value = 'sv'
$ws.Range('J503').Value = 'Statement-opinion'
$ws.Range('I506').Value = 'sd'
$ws.Range('J506').Value = 'Statement-non-opinion'
$ws.Range('I513').Value = 'sd'
$ws.Range('J513').Value = 'Statement-non-opinion'
$ws.Range('I514').Value = 'sd'
$ws.Range('J514').Value = 'Statement-non-opinion'
$ws.Range('I516').Value = '%'
$ws.Range('J516').Value = 'Uninterpretable'
$ws.Range('I521').Value = 'ba'
$ws.Range('J521').Value = 'Appreciation'
$ws.Range('I524').Value = 'sd'
$ws.Range('J524').Value = 'Statement-non-opinion'
